# Correct error in parameters.xlsx
#
# 1. Fix the back-calculated VLDL-flux parameter value on the "parameters"
#    sheet (row 10 / ParamNum 9, kbetaox): 1111.48398681908 -> 1667.2259802286201
#    (dependent shared formulas in columns E and F recalc automatically).
# 2. Add a new "flux summary" worksheet at the end of the workbook with a
#    header row: "Flux from Figure 1" | "Description" | "References".
# 3. Restore page orientation to portrait on the "parameters" sheet.
# 4. Leave the cursor/selection the way the author left it when saving.

$wb = $excel.ActiveWorkbook

$wsParams = $wb.Worksheets.Item("parameters")
$wsUseful = $wb.Worksheets.Item("useful parameters")

# --- 1. Correct the erroneous value -----------------------------------
$wsParams.Range("D10").Value = 1667.2259802286201

# --- 2. Add the new "flux summary" sheet, placed after "useful parameters"
$wsFlux = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsUseful)
$wsFlux.Name = "flux summary"
$wsFlux.Range("A1").Value = "Flux from Figure 1"
$wsFlux.Range("B1").Value = "Description"
$wsFlux.Range("C1").Value = "References"

# --- 3. Page setup on the parameters sheet -----------------------------
$wsParams.PageSetup.Orientation = 1

# --- 4. Final selections / active sheet, matching the saved state ------
$wsUseful.Range("D3").Select()

$wsParams.Activate()
$wsParams.Range("D10").Select()
